$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that is bumped by one day
# (45180 -> 45181) for every data row (rows 2 through 117).
for ($row = 2; $row -le 117; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45180) {
        $cell.Value2 = 45181
    }
}
